$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while preserving it as literal text
# (keeps numeric-looking strings like "245.46" from becoming floats,
# and avoids leaving a lingering per-cell style override).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

# -------- Simple price (column D) refreshes --------
Set-TextValue "D2"  "245.46"
Set-TextValue "D4"  "5.411"
Set-TextValue "D5"  "0.05810"
Set-TextValue "D6"  "3.378"
Set-TextValue "D9"  "0.9966"
Set-TextValue "D10" "0.1424"
Set-TextValue "D11" "0.07503"
Set-TextValue "D13" "0.03054"
Set-TextValue "D14" "4.180"
Set-TextValue "D15" "0.09395"
Set-TextValue "D16" "0.001590"
Set-TextValue "D17" "0.04828"

# -------- Rows 18-24: coin ranking list shifted by one position --------
# Row 18 (was One) -> TigerCash
$ws.Range("B18").Value2 = "TigerCash"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.006238"
$ws.Range("E18").Value2 = "17TigerCashTCH"

# Row 19 (was TigerCash) -> HotbitToken
$ws.Range("B19").Value2 = "HotbitToken"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D19" "0.004092"
$ws.Range("E19").Value2 = "18HotbitTokenHTBWorstin24h"

# Row 20 (was HotbitToken) -> BitKan
$ws.Range("B20").Value2 = "BitKan"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D20" "0.0009983"
$ws.Range("E20").Value2 = "19BitKanKAN"

# Row 21 (was BitKan) -> NitroEx
$ws.Range("B21").Value2 = "NitroEx"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D21" "0.0001501"
$ws.Range("E21").Value2 = "20NitroExNTX"

# Row 22 (was NitroEx) -> LEO
$ws.Range("B22").Value2 = "LEO"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D22" "3.700"
$ws.Range("E22").Value2 = "21LEOLEO"

# Row 23 (was LEO) -> BTSEToken
$ws.Range("B23").Value2 = "BTSEToken"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D23" "2.230"
$ws.Range("E23").Value2 = "22BTSETokenBTSE"

# Row 24 (was BTSEToken) -> One
$ws.Range("B24").Value2 = "One"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D24" "0.0005892"
$ws.Range("E24").Value2 = "23OneONE"

# -------- Remaining simple updates --------
Set-TextValue "D25" "0.3204"
$ws.Range("E27").Value2 = "26UpBotsUBXT"

Set-TextValue "D40" "0.03879"
Set-TextValue "D41" "0.006683"
Set-TextValue "D42" "0.1072"

Set-TextValue "D43" "0.003001"
$ws.Range("E43").Value2 = "42CEJICEJIBestin24h"

Set-TextValue "D44" "0.006458"
Set-TextValue "D48" "0.1459"
